$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect re-pulled / recalculated data
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -8
$ws.Range("F7").Value = -5
$ws.Range("F9").Value = -3
$ws.Range("F12").Value = 5
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 0
$ws.Range("F18").Value = -10
$ws.Range("F19").Value = 3
